# UPDATE: Nama dosen pada select disusun secara ascending A sampai Z
# The combined roster entry "CUT AGUSNIAR, S.T., M.Cs, LIDYA ROSNITA, S.T., M. Kom"
# is split into two separate lecturer rows (CUT AGUSNIAR and LIDYA ROSNITA),
# each carrying their own schedule slots, and a new row is added for LIDYA
# ROSNITA's extra Selasa/Kamis slots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 keeps the Selasa (08:00-10:30) slot and now also picks up the
# Rabu (08:00-10:30) slot that used to live only on row 22; the lecturer
# name becomes just "CUT AGUSNIAR, S.T., M.Cs".
$ws.Range("A21").Value = "CUT AGUSNIAR, S.T., M.Cs"
$ws.Range("D21").Value = "08:00-10:30"

# Row 22 becomes "LIDYA ROSNITA, S.T., M. Kom" and now also carries the
# Selasa (08:00-10:30) slot in addition to the existing Rabu slot.
$ws.Range("A22").Value = "LIDYA ROSNITA, S.T., M. Kom"
$ws.Range("C22").Value = "08:00-10:30"

# New row 27: another LIDYA ROSNITA slot (Selasa 08:00-10:30, Kamis 10:40-13:10).
$ws.Range("A27").Value = "LIDYA ROSNITA, S.T., M. Kom"
$ws.Range("C27").Value = "08:00-10:30"
$ws.Range("E27").Value = "10:40-13:10"

$ws.Range("A27").HorizontalAlignment = -4108
$ws.Range("A27").VerticalAlignment = -4108
$ws.Range("A27").Borders.Item(7).LineStyle = 1
$ws.Range("A27").Borders.Item(10).LineStyle = 1

# Scroll/zoom/selection state left behind by the editor when the change was made.
$ws.Application.ActiveWindow.Zoom = 115
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("A22").Select()
